# Auto-generated Excel COM-interop edit script
# Applies the cryptos.xlsx price/volume update described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.708.69'
$ws.Range('E2').Value = '  +0.03%  '
$ws.Range('D3').Value = '1.636.34'
$ws.Range('E3').Value = '  -0.61%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '''217.35'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.71%  '
$ws.Range('D6').Value = '''0.502'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.77%  '
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('D8').Value = '''0.250'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.77%  '
$ws.Range('D9').Value = '''0.0622'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.79%  '
$ws.Range('D10').Value = '''19.06'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.44%  '
$ws.Range('D11').Value = '''0.0843'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.02%  '
$ws.Range('D12').Value = '1.866.67'
$ws.Range('E12').Value = '  -0.56%  '
$ws.Range('D13').Value = '1.635.91'
$ws.Range('E13').Value = '  -0.45%  '
$ws.Range('D14').Value = '''4.13'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.08%  '
$ws.Range('D15').Value = '''0.525'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.26%  '
$ws.Range('D16').Value = '''64.30'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.29%  '
$ws.Range('D17').Value = '26.707.01'
$ws.Range('E17').Value = '  -0.04%  '
$ws.Range('D18').Value = '0.0₃0727'
$ws.Range('E18').Value = '  -2.26%  '
$ws.Range('B19').Value = 'Dai'
$ws.Range('C19').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D19').Value = '''1.01'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.14%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').Value = '''210.62'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.08%  '
$ws.Range('D21').Value = '''4.34'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.38%  '
$ws.Range('D22').Value = '''2.37'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.29%  '
$ws.Range('D23').Value = '''6.18'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.45%  '
$ws.Range('D24').Value = '''9.23'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.86%  '
$ws.Range('D25').Value = '''145.45'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.24%  '
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('E27').Value = '  -1.84%  '
$ws.Range('D28').Value = '''7.09'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.44%  '
$ws.Range('D29').Value = '''15.57'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.02%  '
$ws.Range('D30').Value = '''0.0504'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.02%  '
$ws.Range('D31').Value = '''1.18'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.65%  '
$ws.Range('D32').Value = '''3.35'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.21%  '
$ws.Range('D33').Value = '''2.97'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.14%  '
$ws.Range('D34').Value = '1.273.45'
$ws.Range('E34').Value = '  -0.26%  '
$ws.Range('E35').Value = '  -0.98%  '
$ws.Range('E36').Value = '  +0.35%  '
$ws.Range('E37').Value = '  -2.09%  '
$ws.Range('D38').Value = '''0.530'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.69%  '
$ws.Range('D39').Value = '''0.807'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.71%  '
$ws.Range('E40').Value = '  +0.08%  '
$ws.Range('E41').Value = '  -1.53%  '
$ws.Range('E42').Value = '  -1.91%  '
$ws.Range('D43').Value = '1.777.05'
$ws.Range('E43').Value = '  -0.51%  '
$ws.Range('E44').Value = '  -3.54%  '
$ws.Range('D45').Value = '''60.70'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.48%  '
$ws.Range('D46').Value = '''91.08'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.65%  '
$ws.Range('D47').Value = '''1.57'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.96%  '
$ws.Range('D48').Value = '''0.0519'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.67%  '
$ws.Range('D49').Value = '''7.58'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.46%  '
$ws.Range('D50').Value = '''0.0958'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.82%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').Value = '''0.407'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.22%  '
